# Update "想去人数" (number of people interested) figures in column F
# across the four sheets of the workbook, matching the refreshed scrape
# output committed as "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value  = 1564
$ws1.Cells.Item(4, 6).Value  = 816
$ws1.Cells.Item(5, 6).Value  = 225
$ws1.Cells.Item(6, 6).Value  = 60
$ws1.Cells.Item(7, 6).Value  = 1122
$ws1.Cells.Item(8, 6).Value  = 716
$ws1.Cells.Item(9, 6).Value  = 771
$ws1.Cells.Item(10, 6).Value = 1389
$ws1.Cells.Item(12, 6).Value = 1023
$ws1.Cells.Item(15, 6).Value = 186
$ws1.Cells.Item(17, 6).Value = 435
$ws1.Cells.Item(18, 6).Value = 10
$ws1.Cells.Item(20, 6).Value = 291
$ws1.Cells.Item(21, 6).Value = 540
$ws1.Cells.Item(22, 6).Value = 558
$ws1.Cells.Item(23, 6).Value = 750
$ws1.Cells.Item(24, 6).Value = 235
$ws1.Cells.Item(25, 6).Value = 169

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(3, 6).Value  = 991
$ws2.Cells.Item(5, 6).Value  = 252
$ws2.Cells.Item(7, 6).Value  = 139
$ws2.Cells.Item(8, 6).Value  = 63
$ws2.Cells.Item(10, 6).Value = 76

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 209

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value  = 209
$ws4.Cells.Item(4, 6).Value  = 1564
$ws4.Cells.Item(6, 6).Value  = 816
$ws4.Cells.Item(7, 6).Value  = 225
$ws4.Cells.Item(8, 6).Value  = 991
$ws4.Cells.Item(9, 6).Value  = 60
$ws4.Cells.Item(10, 6).Value = 1122
$ws4.Cells.Item(11, 6).Value = 716
$ws4.Cells.Item(12, 6).Value = 771
$ws4.Cells.Item(13, 6).Value = 1389
$ws4.Cells.Item(15, 6).Value = 1023
$ws4.Cells.Item(18, 6).Value = 186
$ws4.Cells.Item(20, 6).Value = 435
$ws4.Cells.Item(21, 6).Value = 10
$ws4.Cells.Item(23, 6).Value = 252
$ws4.Cells.Item(25, 6).Value = 291
$ws4.Cells.Item(27, 6).Value = 139
$ws4.Cells.Item(28, 6).Value = 139
$ws4.Cells.Item(29, 6).Value = 540
$ws4.Cells.Item(30, 6).Value = 558
$ws4.Cells.Item(31, 6).Value = 750
$ws4.Cells.Item(32, 6).Value = 235
$ws4.Cells.Item(33, 6).Value = 63
$ws4.Cells.Item(34, 6).Value = 169
$ws4.Cells.Item(36, 6).Value = 76
$ws4.Cells.Item(37, 6).Value = 76
